# "obs solo con lp" - fill in the observed data for the LP (Probing) table
# on sheet "Datos Lab7": columns B (Consumo de Datos [kB]) and C (Tiempo de
# Ejecucion Real @LP [ms]) for rows 3:6 (Factor de Carga 0.1/0.5/0.7/0.9).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Datos Lab7")

$ws.Range("B3").Value = 2210470.146
$ws.Range("C3").Value = 99410.392

$ws.Range("B4").Value = 1820936.711
$ws.Range("C4").Value = 83507.486

$ws.Range("B5").Value = 1712145.346
$ws.Range("C5").Value = 88886.338

$ws.Range("B6").Value = 1593599.386
$ws.Range("C6").Value = 306333.314
